$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("FFOptionMappingFinal")
$ws3.Rows("602:604").Insert()
$ws3.Range("C602").NumberFormat = "@"
$ws3.Range("C602").Value = "68"
$ws3.Range("C601").Copy()
$ws3.Range("C602").PasteSpecial(-4122)
